# Translations update pulled in ahead of an upstream merge.
# Inserts a new "Sheet1" (string-token translation table) between the
# existing "settings" and "choices" tabs, and makes it the active tab.

$wb = $excel.ActiveWorkbook

# Insert the new worksheet right after "settings" (i.e. before "choices").
$settings = $wb.Worksheets.Item("settings")
$translations = $wb.Worksheets.Add($null, $settings)

# Header row.
$translations.Range("A1").Value = "string_token"
$translations.Range("B1").Value = "text.default"
$translations.Range("C1").Value = "text.spanish"

# Token rows.
$translations.Range("A2").Value = "delivery_id"
$translations.Range("B2").Value = "Delivery ID"

$translations.Range("A3").Value = "date_time"
$translations.Range("B3").Value = "Time of Delivery"

$translations.Range("A4").Value = "delivery_site"
$translations.Range("B4").Value = "Delivery Site"

$translations.Range("A5").Value = "distributor"
$translations.Range("B5").Value = "Distributor"

$translations.Range("A6").Value = "scanned_item_pack"
$translations.Range("B6").Value = "Scanned Item Pack"

# Match the selection state captured in the target workbook and make this
# new sheet the active / visibly-selected tab.
$translations.Range("C3").Select()
$translations.Activate()
